$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number that was updated
# from 45178 (2023-09-09) to 45179 (2023-09-10) for every data row
# (rows 2 through 375).
$ws.Range("C2:C375").Value = 45179
